# Updates cell values in the cryptocurrency price table (Price/Volume(1h)
# columns) to reflect the latest scrape, per the GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage (leading apostrophe,
# exactly like typing `'69.55` into Excel) so that price strings which look
# like plain numbers ("0.520", "7.57", ...) are not silently re-interpreted
# as numeric values (which would drop meaningful trailing zeros). Values that
# already contain multiple separators (e.g. "69.566.51") or formatting that
# prevents numeric parsing do not need this and are assigned directly.
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
}

$ws.Range("D2").Value = '69.566.51'
$ws.Range("E2").Value = '  -1.46%  '
$ws.Range("D3").Value = '2.515.85'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("E4").Value = '  +0.04%  '
Set-TextValue $ws.Range("D5") '573.25'
$ws.Range("E5").Value = '  -0.76%  '
Set-TextValue $ws.Range("D6") '166.16'
$ws.Range("E6").Value = '  -1.95%  '
$ws.Range("E7").Value = '  -0.07%  '
Set-TextValue $ws.Range("D8") '0.514'
$ws.Range("E8").Value = '  +0.41%  '
$ws.Range("D9").Value = '2.514.34'
$ws.Range("E9").Value = '  -0.37%  '
$ws.Range("E10").Value = '  -1.87%  '
$ws.Range("E11").Value = '  -0.57%  '
Set-TextValue $ws.Range("D12") '0.357'
$ws.Range("E12").Value = '  +3.94%  '
$ws.Range("E13").Value = '  +1.07%  '
$ws.Range("D14").Value = '2.974.75'
$ws.Range("E14").Value = '  -0.38%  '
$ws.Range("D15").Value = '69.409.92'
$ws.Range("E15").Value = '  -1.54%  '
$ws.Range("E16").Value = '  -2.24%  '
Set-TextValue $ws.Range("D17") '24.81'
$ws.Range("E17").Value = '  -1.51%  '
$ws.Range("D18").Value = '2.517.87'
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("E19").Value = '  -1.74%  '
Set-TextValue $ws.Range("D20") '7.57'
$ws.Range("E20").Value = '  -0.88%  '
Set-TextValue $ws.Range("D21") '350.26'
$ws.Range("E21").Value = '  -2.82%  '
$ws.Range("E22").Value = '  -1.47%  '
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("E24").Value = '  -0.04%  '
Set-TextValue $ws.Range("D25") '70.33'
$ws.Range("E25").Value = '  +1.25%  '
Set-TextValue $ws.Range("D26") '3.92'
$ws.Range("E26").Value = '  -3.62%  '
Set-TextValue $ws.Range("D27") '8.91'
$ws.Range("E27").Value = '  -2.76%  '
$ws.Range("D28").Value = '2.644.33'
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("E29").Value = '  +0.72%  '
$ws.Range("D30").Value = '0.0₃0891'
$ws.Range("E30").Value = '  -2.95%  '
Set-TextValue $ws.Range("D31") '7.80'
$ws.Range("E31").Value = '  -0.74%  '
Set-TextValue $ws.Range("D32") '463.38'
$ws.Range("E32").Value = '  -4.62%  '
Set-TextValue $ws.Range("D33") '1.24'
$ws.Range("E33").Value = '  -4.87%  '
$ws.Range("E34").Value = '  -2.13%  '
$ws.Range("E35").Value = '  +0.01%  '
Set-TextValue $ws.Range("D36") '0.117'
$ws.Range("E36").Value = '  +1.43%  '
Set-TextValue $ws.Range("D37") '157.33'
$ws.Range("E37").Value = '  +0.80%  '
$ws.Range("E38").Value = '  +0.98%  '
Set-TextValue $ws.Range("D39") '18.48'
$ws.Range("E39").Value = '  -0.97%  '
$ws.Range("E41").Value = '  -1.04%  '
$ws.Range("E42").Value = '  -1.89%  '
$ws.Range("E43").Value = '  -2.93%  '
$ws.Range("E46").Value = '  -7.35%  '
$ws.Range("E47").Value = '  -1.47%  '
Set-TextValue $ws.Range("D48") '0.520'
$ws.Range("E48").Value = '  -1.69%  '
Set-TextValue $ws.Range("D49") '3.47'
$ws.Range("E49").Value = '  -2.24%  '
Set-TextValue $ws.Range("D50") '0.0733'
$ws.Range("E50").Value = '  +0.22%  '
Set-TextValue $ws.Range("D51") '0.578'
$ws.Range("E51").Value = '  -3.76%  '
